# Update "想去人数" (want-to-go count) figures across the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2521
$ws1.Range("F8").Value  = 926
$ws1.Range("F11").Value = 1944
$ws1.Range("F12").Value = 673
$ws1.Range("F21").Value = 6

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 231

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 578
$ws3.Range("F4").Value = 561

# --- Sheet "全部类型" (all types, merged view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 578
$ws4.Range("F5").Value  = 561
$ws4.Range("F15").Value = 2521
$ws4.Range("F16").Value = 926
$ws4.Range("F22").Value = 231
$ws4.Range("F23").Value = 231
$ws4.Range("F24").Value = 1944
$ws4.Range("F25").Value = 673
$ws4.Range("F38").Value = 6
